# Update Name of Algo
# Apply corrected numeric values to the result_data_KNN sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -12.934
$ws.Range("E6").Value = 12.967
$ws.Range("E7").Value = 13.051
$ws.Range("C8").Value = -12.666
$ws.Range("E8").Value = 13.084
$ws.Range("B12").Value = 5.417
$ws.Range("C12").Value = -13.073
$ws.Range("C14").Value = -11.675
$ws.Range("E19").Value = 12.614
$ws.Range("E21").Value = 12.93
$ws.Range("C22").Value = -12.929
$ws.Range("E24").Value = 12.855
